# Powerpoint writer: consolidate text runs when possible.
# This applies the same consolidation the golden test expects: several
# shapes had their text split across many <a:r> runs (one run per word,
# with blank-space runs in between). We rewrite each TextRange so the
# whole paragraph collapses into a single run.
#
# Because the COM text-range setter treats assignments that share a
# text prefix with the existing (split) text as an incremental edit
# (only appending/inserting the differing suffix, leaving the old runs
# untouched), we first assign a disjoint placeholder string and then
# assign the real text. That forces a full-run replacement, yielding a
# single <a:r><a:rPr/><a:t>...</a:t></a:r> run.

$p = $ppt.ActivePresentation

function Set-ConsolidatedText($shape, [string]$text) {
    $tr = $shape.TextFrame.TextRange
    $tr.Text = "##RUN-CONSOLIDATION-PLACEHOLDER##"
    $tr.Text = $text
}

# Slide 1: "Section Header (with background image)" (Title)
$s1 = $p.Slides.Item(1)
Set-ConsolidatedText $s1.Shapes.Item(1) "Section Header (with background image)"

# Slide 2: "Slide 1" (Title)
$s2 = $p.Slides.Item(2)
Set-ConsolidatedText $s2.Shapes.Item(1) "Slide 1"

# Slide 3: "Slide 2" (Title)
$s3 = $p.Slides.Item(3)
Set-ConsolidatedText $s3.Shapes.Item(1) "Slide 2"

# Slide 4: "Slide 3" (Title)
$s4 = $p.Slides.Item(4)
Set-ConsolidatedText $s4.Shapes.Item(1) "Slide 3"

# Slide 5: "Slide 4" (Title) and "An image" (TextBox 3)
$s5 = $p.Slides.Item(5)
Set-ConsolidatedText $s5.Shapes.Item(1) "Slide 4"
Set-ConsolidatedText $s5.Shapes.Item(4) "An image"

# Slide 6 notes: "Blank slides can have background images." (Notes Placeholder 2)
$s6 = $p.Slides.Item(6)
$notesShape = $s6.NotesPage.Shapes.Item(2)
Set-ConsolidatedText $notesShape "Blank slides can have background images."
